$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Headers")
$ws.Range("A23").Value = "TEST"
